$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- Shape "Freeform 2" (COM index 2) : text "node* next" -> cleared ---
$shNext = $s.Shapes.Item(2)
$shNext.TextFrame.TextRange.Text = ""

# --- Shape "TextBox 4" (COM index 4) : reposition + split runs ---
$shFields = $s.Shapes.Item(4)

# Reposition (EMU -> points, 1 pt = 12700 EMU)
$shFields.Left = 1018441 / 12700
$shFields.Top  = 3064398 / 12700

$tr = $shFields.TextFrame.TextRange

# Paragraph 1 (chars 1-12): "node* p_next" -> runs "node* " + "p_next"
$run1a = $tr.Characters(1, 6)
$run1a.Text = "node* "
$run1b = $tr.Characters(7, 6)
$run1b.Text = "p_next"

# Paragraph 2 (chars 14-31, after the para-break at char 13): "const char* p_data"
# -> runs "const char* " + "p_data"
$run2a = $tr.Characters(14, 12)
$run2a.Text = "const char* "
$run2b = $tr.Characters(26, 6)
$run2b.Text = "p_data"

Write-Host "Updated slide 6: cleared 'node* next' label and split node*/p_next, const char*/p_data runs."
